$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 9035.083000000001
$ws.Cells.Item(9, 9).Value = 12741.125
$ws.Cells.Item(9, 10).Value = 1623
$ws.Cells.Item(9, 11).Value = 12741.125
$ws.Cells.Item(9, 12).Value = 1623
$ws.Cells.Item(9, 13).Value = -12572.125
$ws.Cells.Item(9, 14).Value = -1961
$ws.Cells.Item(19, 8).Value = 3436.1177
$ws.Cells.Item(19, 9).Value = 1838.4445
$ws.Cells.Item(19, 10).Value = 5233.5
$ws.Cells.Item(19, 11).Value = 1838.4445
$ws.Cells.Item(19, 12).Value = 5233.5
$ws.Cells.Item(19, 13).Value = -1663.4445
$ws.Cells.Item(19, 14).Value = -5583.5
$ws.Cells.Item(43, 8).Value = 717118.3
$ws.Cells.Item(43, 9).Value = 3206.8572
$ws.Cells.Item(43, 10).Value = 1431029.8
$ws.Cells.Item(43, 11).Value = 3206.8572
$ws.Cells.Item(43, 12).Value = 1431029.8
$ws.Cells.Item(43, 13).Value = -3137.8572
$ws.Cells.Item(43, 14).Value = -1431167.8
$ws.Cells.Item(47, 8).Value = 22983.334
$ws.Cells.Item(47, 9).Value = 24475
$ws.Cells.Item(47, 11).Value = 24475
$ws.Cells.Item(47, 13).Value = -23503
$ws.Cells.Item(88, 8).Value = 1841.625
$ws.Cells.Item(88, 9).Value = 2121.5
$ws.Cells.Item(88, 11).Value = 2121.5
$ws.Cells.Item(88, 13).Value = -1715.5
$ws.Cells.Item(91, 8).Value = 1841.625
$ws.Cells.Item(91, 9).Value = 2121.5
$ws.Cells.Item(91, 11).Value = 2121.5
$ws.Cells.Item(91, 13).Value = -717.5
$ws.Cells.Item(121, 8).Value = 2465.2856
$ws.Cells.Item(121, 10).Value = 2465.2856
$ws.Cells.Item(121, 12).Value = 7395.8568
$ws.Cells.Item(121, 14).Value = -10889.8568
$ws.Cells.Item(132, 8).Value = 2275.1353
$ws.Cells.Item(132, 9).Value = 2278.3438
$ws.Cells.Item(132, 11).Value = 6835.0314
$ws.Cells.Item(132, 13).Value = -4305.0314
$ws.Cells.Item(137, 8).Value = 3209.3333
$ws.Cells.Item(137, 9).Value = 3662.4
$ws.Cells.Item(137, 10).Value = 2076.6667
$ws.Cells.Item(137, 11).Value = 10987.2
$ws.Cells.Item(137, 12).Value = 6230.000100000001
$ws.Cells.Item(137, 13).Value = -8437.200000000001
$ws.Cells.Item(137, 14).Value = -11330.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 582.1818
$ws.Cells.Item(5, 9).Value = 540.2
$ws.Cells.Item(5, 11).Value = 540.2
$ws.Cells.Item(5, 13).Value = -428.2
$ws.Cells.Item(61, 8).Value = 27029716
$ws.Cells.Item(61, 9).Value = 34484036
$ws.Cells.Item(61, 10).Value = 7814.25
$ws.Cells.Item(61, 11).Value = 34484036
$ws.Cells.Item(61, 12).Value = 7814.25
$ws.Cells.Item(61, 13).Value = -34483824
$ws.Cells.Item(61, 14).Value = -8238.25
$ws.Cells.Item(97, 8).Value = 1181.125
$ws.Cells.Item(97, 9).Value = 564.1429000000001
$ws.Cells.Item(97, 10).Value = 5500
$ws.Cells.Item(97, 11).Value = 564.1429000000001
$ws.Cells.Item(97, 12).Value = 5500
$ws.Cells.Item(97, 13).Value = -68.14290000000005
$ws.Cells.Item(97, 14).Value = -6492
$ws.Cells.Item(124, 8).Value = 44744.285
$ws.Cells.Item(124, 10).Value = 44744.285
$ws.Cells.Item(124, 12).Value = 44744.285
$ws.Cells.Item(124, 14).Value = -54564.285
$ws.Cells.Item(136, 8).Value = 27029716
$ws.Cells.Item(136, 9).Value = 34484036
$ws.Cells.Item(136, 10).Value = 7814.25
$ws.Cells.Item(136, 11).Value = 103452108
$ws.Cells.Item(136, 12).Value = 23442.75
$ws.Cells.Item(136, 13).Value = -103449558
$ws.Cells.Item(136, 14).Value = -28542.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 582.1818
$ws.Cells.Item(4, 9).Value = 540.2
$ws.Cells.Item(4, 11).Value = 540.2
$ws.Cells.Item(4, 13).Value = -425.2
$ws.Cells.Item(9, 8).Value = 50000
$ws.Cells.Item(9, 10).Value = 50000
$ws.Cells.Item(9, 12).Value = 50000
$ws.Cells.Item(9, 14).Value = -50336
$ws.Cells.Item(20, 8).Value = 2228.2778
$ws.Cells.Item(20, 9).Value = 2201.5386
$ws.Cells.Item(20, 11).Value = 2201.5386
$ws.Cells.Item(20, 13).Value = -1954.5386
$ws.Cells.Item(74, 8).Value = 39911.6
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 39911.6
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 39911.6
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(74, 14).Value = -41783.6
$ws.Cells.Item(75, 8).Value = 22470.637
$ws.Cells.Item(75, 9).Value = 4461
$ws.Cells.Item(75, 10).Value = 53987.5
$ws.Cells.Item(75, 11).Value = 4461
$ws.Cells.Item(75, 12).Value = 53987.5
$ws.Cells.Item(75, 13).Value = -3525
$ws.Cells.Item(75, 14).Value = -55859.5
$ws.Cells.Item(77, 8).Value = 39911.6
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 39911.6
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 119734.8
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(77, 14).Value = -129094.8
$ws.Cells.Item(78, 8).Value = 22470.637
$ws.Cells.Item(78, 9).Value = 4461
$ws.Cells.Item(78, 10).Value = 53987.5
$ws.Cells.Item(78, 11).Value = 13383
$ws.Cells.Item(78, 12).Value = 161962.5
$ws.Cells.Item(78, 13).Value = -8703
$ws.Cells.Item(78, 14).Value = -171322.5
$ws.Cells.Item(105, 8).Value = 15491.889
$ws.Cells.Item(105, 9).Value = 19489.5
$ws.Cells.Item(105, 10).Value = 7496.6665
$ws.Cells.Item(105, 11).Value = 19489.5
$ws.Cells.Item(105, 12).Value = 7496.6665
$ws.Cells.Item(105, 13).Value = -17742.5
$ws.Cells.Item(105, 14).Value = -10990.6665
$ws.Cells.Item(134, 8).Value = 4755.9287
$ws.Cells.Item(134, 9).Value = 4755.9287
$ws.Cells.Item(134, 11).Value = 14267.7861
$ws.Cells.Item(134, 13).Value = -11732.7861

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 9).Value = 236.625
$ws.Cells.Item(7, 10).Value = 438.375
$ws.Cells.Item(7, 11).Value = 236.625
$ws.Cells.Item(7, 12).Value = 438.375
$ws.Cells.Item(7, 13).Value = -123.625
$ws.Cells.Item(7, 14).Value = -664.375
$ws.Cells.Item(31, 8).Value = 23813344
$ws.Cells.Item(31, 9).Value = 3293.8235
$ws.Cells.Item(31, 10).Value = 125006060
$ws.Cells.Item(31, 11).Value = 3293.8235
$ws.Cells.Item(31, 12).Value = 125006060
$ws.Cells.Item(31, 13).Value = -2998.8235
$ws.Cells.Item(31, 14).Value = -125006650
$ws.Cells.Item(34, 8).Value = 23813344
$ws.Cells.Item(34, 9).Value = 3293.8235
$ws.Cells.Item(34, 10).Value = 125006060
$ws.Cells.Item(34, 11).Value = 3293.8235
$ws.Cells.Item(34, 12).Value = 125006060
$ws.Cells.Item(34, 13).Value = -3091.8235
$ws.Cells.Item(34, 14).Value = -125006464
$ws.Cells.Item(122, 8).Value = 1200.2084
$ws.Cells.Item(122, 9).Value = 1262.5333
$ws.Cells.Item(122, 11).Value = 3787.5999
$ws.Cells.Item(122, 13).Value = -1337.5999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 10150.538
$ws.Cells.Item(133, 10).Value = 19496.5
$ws.Cells.Item(133, 12).Value = 58489.5
$ws.Cells.Item(133, 14).Value = -68609.5
$ws.Cells.Item(136, 8).Value = 5498.4287
$ws.Cells.Item(136, 9).Value = 4997.25
$ws.Cells.Item(136, 11).Value = 14991.75
$ws.Cells.Item(136, 13).Value = -9891.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 1250476.1
$ws.Cells.Item(14, 9).Value = 1428829.9
$ws.Cells.Item(14, 10).Value = 2000
$ws.Cells.Item(14, 11).Value = 1428829.9
$ws.Cells.Item(14, 12).Value = 2000
$ws.Cells.Item(14, 13).Value = -1428661.9
$ws.Cells.Item(14, 14).Value = -2336
$ws.Cells.Item(57, 8).Value = 19000
$ws.Cells.Item(57, 10).Value = 20000
$ws.Cells.Item(57, 12).Value = 20000
$ws.Cells.Item(57, 14).Value = -21640
$ws.Cells.Item(70, 8).Value = 4493
$ws.Cells.Item(70, 9).Value = 4366.25
$ws.Cells.Item(70, 11).Value = 4366.25
$ws.Cells.Item(70, 13).Value = -4096.25
$ws.Cells.Item(73, 8).Value = 4493
$ws.Cells.Item(73, 9).Value = 4366.25
$ws.Cells.Item(73, 11).Value = 4366.25
$ws.Cells.Item(73, 13).Value = -3430.25
$ws.Cells.Item(107, 8).Value = 335.92593
$ws.Cells.Item(107, 9).Value = 421.6875
$ws.Cells.Item(107, 10).Value = 211.18182
$ws.Cells.Item(107, 11).Value = 421.6875
$ws.Cells.Item(107, 12).Value = 211.18182
$ws.Cells.Item(107, 13).Value = 1498.3125
$ws.Cells.Item(107, 14).Value = -4051.18182
$ws.Cells.Item(132, 8).Value = 3179.2917
$ws.Cells.Item(132, 9).Value = 3264.6667
$ws.Cells.Item(132, 10).Value = 2581.6667
$ws.Cells.Item(132, 11).Value = 9794.000100000001
$ws.Cells.Item(132, 12).Value = 7745.000100000001
$ws.Cells.Item(132, 13).Value = -7264.000100000001
$ws.Cells.Item(132, 14).Value = -12805.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2241.5715
$ws.Cells.Item(22, 9).Value = 1444.3889
$ws.Cells.Item(22, 10).Value = 2839.4583
$ws.Cells.Item(22, 11).Value = 1444.3889
$ws.Cells.Item(22, 12).Value = 2839.4583
$ws.Cells.Item(22, 13).Value = -1149.3889
$ws.Cells.Item(22, 14).Value = -3429.4583
$ws.Cells.Item(27, 8).Value = 2241.5715
$ws.Cells.Item(27, 9).Value = 1444.3889
$ws.Cells.Item(27, 10).Value = 2839.4583
$ws.Cells.Item(27, 11).Value = 1444.3889
$ws.Cells.Item(27, 12).Value = 2839.4583
$ws.Cells.Item(27, 13).Value = -1337.3889
$ws.Cells.Item(27, 14).Value = -3053.4583

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1561.32
$ws.Cells.Item(136, 9).Value = 1420.375
$ws.Cells.Item(136, 10).Value = 4944
$ws.Cells.Item(136, 11).Value = 4261.125
$ws.Cells.Item(136, 12).Value = 14832
$ws.Cells.Item(136, 13).Value = -1711.125
$ws.Cells.Item(136, 14).Value = -19932
$ws.Cells.Item(137, 8).Value = 135000
$ws.Cells.Item(137, 10).Value = 135000
$ws.Cells.Item(137, 12).Value = 135000
$ws.Cells.Item(137, 14).Value = -145200
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 150000
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 100788
$ws.Cells.Item(141, 9).Value = 77980
$ws.Cells.Item(141, 10).Value = 135000
$ws.Cells.Item(141, 11).Value = 77980
$ws.Cells.Item(141, 12).Value = 135000
$ws.Cells.Item(141, 13).Value = -72800
$ws.Cells.Item(141, 14).Value = -145360
